# Domänenmodell_text.docx - textuelle Beschreibung Domänenmodell verbessert
#
# Strategy: for each paragraph that needs restructuring (run splits, bookmark
# relocation, text changes), replace the paragraph's full range (including its
# end-of-paragraph mark) with freshly authored WordprocessingML via
# Range.InsertXML, wrapped in the standard pkg:package envelope Word expects.
# This gives exact control over run boundaries and bookmark placement, which
# plain Find/Replace or InsertAfter cannot guarantee (they tend to merge text
# into existing runs instead of creating new ones).

$d = $word.ActiveDocument

function Set-ParagraphXml($paragraphIndex, [string]$innerXml) {
    $para = $d.Paragraphs($paragraphIndex).Range
    $fullXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>$innerXml</w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
    $para.InsertXML($fullXml)
}

# 1) "Domänenmodell:" -> "Domänenmodell" + bookmark "_GoBack" + ":"
Set-ParagraphXml 1 '<w:r><w:t>Domänenmodell</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>:</w:t></w:r>'

# 2) "Jeder Bewohner hat einen Namen und ein Passwort" -> add trailing "." run
Set-ParagraphXml 3 '<w:r><w:t>Jeder Bewohner hat einen Namen und ein Passwort</w:t></w:r><w:r><w:t>.</w:t></w:r>'

# 3) "Dem Bewohner werden Aufgaben zugeteilt" -> add trailing "." run
Set-ParagraphXml 5 '<w:r><w:t>Dem Bewohner werden Aufgaben zugeteilt</w:t></w:r><w:r><w:t>.</w:t></w:r>'

# 4) "Die Elemente der Aufgabenliste..." -> drop the _GoBack bookmark (it moved to paragraph 1)
Set-ParagraphXml 8 '<w:r><w:t>Die Elemente der Aufgabenliste ergeben sich aus der Raumbeschreibung.</w:t></w:r>'

# 5) "Der Bewohner verwendet eine Einkaufsliste." -> reworded, split across 4 runs
Set-ParagraphXml 12 '<w:r><w:t>Alle Bewohner verwenden</w:t></w:r><w:r><w:t xml:space="preserve"> eine </w:t></w:r><w:r><w:t xml:space="preserve">gemeinsame </w:t></w:r><w:r><w:t>Einkaufsliste.</w:t></w:r>'

# 6) "Einem Bewohner wird ein Privatkalender zugeteilt." -> expanded sentence, split across 2 runs
Set-ParagraphXml 14 '<w:r><w:t>Einem Bewohner wi</w:t></w:r><w:r><w:t xml:space="preserve">rd ein Privatkalender zugeteilt, in den er seine nur für sich sichtbaren Termine eintragen kann.  </w:t></w:r>'

# 7) "Der Gemeinschaftsplan ist für beliebig viele Bewohner zugänglich. " -> replaced sentence
Set-ParagraphXml 15 '<w:r><w:t xml:space="preserve">In den Gemeinschaftsplan können alle Bewohner gemeinsame Termine eintragen. Diese sind dann für alle Bewohner im Gemeinschaftsplan einsehbar. </w:t></w:r>'

# 8) "Die Kalender bestehen aus Terminen." -> reworded, split across 3 runs
Set-ParagraphXml 17 '<w:r><w:t xml:space="preserve">Die Kalender </w:t></w:r><w:r><w:t>enthalten</w:t></w:r><w:r><w:t xml:space="preserve"> Terminen.</w:t></w:r>'

Write-Host "Paragraphs after edit: $($d.Paragraphs.Count)"
